# Feature: output pine scripts to draw trading info.
#
# The "comment"/"done" flag (shared string "Yes") that marked row 87 as
# finished is moved down one row (the underlying parameter rows shifted),
# and a brand-new parameter-sweep row (93) is appended with its own
# "Yes" flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Yes" flag that used to sit on row 87 no longer belongs there ...
$ws.Range("B87").ClearContents()

# ... it now belongs on rows 89-92.
$ws.Range("B89").Value = "Yes"
$ws.Range("B90").Value = "Yes"
$ws.Range("B91").Value = "Yes"
$ws.Range("B92").Value = "Yes"

# Append a new parameter row (93), cloning the date formatting used by
# the row above it so the new dates keep the short-date number format.
$ws.Range("F92:G92").Copy() | Out-Null
$ws.Range("F93:G93").PasteSpecial(-4122) | Out-Null

$ws.Range("A93").Value = 92
$ws.Range("B93").Value = "Yes"
$ws.Range("C93").Value = 200000
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = "ETH"
$ws.Range("F93").Value = (Get-Date -Year 2021 -Month 1 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("G93").Value = (Get-Date -Year 2021 -Month 2 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Range("H93").Value = 0.615946
$ws.Range("I93").Value = 0.999249
$ws.Range("J93").Value = 0.008156
$ws.Range("K93").Value = 0.740697

# Match the author's final cursor position.
$ws.Range("F96").Select() | Out-Null
